$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Deadline: August 24th" -> "Deadline: August 30th"
# ------------------------------------------------------------------
$pDeadline = $d.Paragraphs.Item(5)
$pDeadline.Range.Find.Execute("24", $true, $false, $false, $false, $false, `
    $true, 1, $false, "30", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "August 22-23rd ~ Finish GUI " -> "August 28-29th ~ Finish GUI "
#    (with the "rd" -> "th" ordinal, and a "_GoBack" bookmark now sitting
#    right after the "~" - this also relocates the single document-wide
#    "_GoBack" bookmark away from its old position automatically)
# ------------------------------------------------------------------
$pGui = $d.Paragraphs.Item(6)
$pGui.Range.Find.Execute("August 22-23", $true, $false, $false, $false, $false, `
    $true, 1, $false, "August 28-29", 2) | Out-Null
$pGui.Range.Find.Execute("rd", $true, $false, $false, $false, $false, `
    $true, 1, $false, "th", 2) | Out-Null

# find the "~" inside this paragraph and plant the bookmark right after it
$tildeRange = $d.Range($pGui.Range.Start, $pGui.Range.End)
$tildeRange.Find.Execute("~", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $d.Range($tildeRange.End, $tildeRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPos) | Out-Null

# ------------------------------------------------------------------
# 3) "August 16-20th" -> "August 22-26th"
# ------------------------------------------------------------------
$pSms = $d.Paragraphs.Item(7)
$pSms.Range.Find.Execute("August 16-20", $true, $false, $false, $false, $false, `
    $true, 1, $false, "August 22-26", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "August 9-10th" -> "August 15-16th"  (Start SMS part line)
# ------------------------------------------------------------------
$pStartSms = $d.Paragraphs.Item(8)
$pStartSms.Range.Find.Execute("August 9-10", $true, $false, $false, $false, $false, `
    $true, 1, $false, "August 15-16", 2) | Out-Null

# ------------------------------------------------------------------
# 5) "August 7-8th" -> "August 12-14th"
# ------------------------------------------------------------------
$pText = $d.Paragraphs.Item(9)
$pText.Range.Find.Execute("August 7-8", $true, $false, $false, $false, $false, `
    $true, 1, $false, "August 12-14", 2) | Out-Null

# ------------------------------------------------------------------
# 6) New bullet after the "Optional Closing Message" line:
#    "Creating the Forecast Body - August 9-10th"
# ------------------------------------------------------------------
$pClosing = $d.Paragraphs.Item(12)
$pClosing.Range.InsertParagraphAfter() | Out-Null

$pNew = $d.Paragraphs.Item(13)
$pNew.Range.InsertBefore("Creating the Forecast Body " + [char]0x2013 + " August 9-10") | Out-Null

$ordEnd = $pNew.Range.End - 1
$ordRange = $d.Range($ordEnd, $ordEnd)
$ordRange.InsertAfter("th") | Out-Null
$ordRange.Font.Superscript = $true
